# Update "想去人数" (interested-in count) values in column F, rows 3-8,
# on both the "展览" and "全部类型" sheets to reflect the new scrape.

$wb = $excel.ActiveWorkbook

$updates = @{
    3 = 1371
    4 = 10
    5 = 8
    6 = 4
    7 = 37
    8 = 190
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
